$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'286.41"
$ws.Range("D3").Value = "'21.17"
$ws.Range("D4").Value = "'6.451"
$ws.Range("D5").Value = "'0.06367"
$ws.Range("D7").Value = "'1.557"
$ws.Range("D8").Value = "'6.579"
$ws.Range("D9").Value = "'0.8237"
$ws.Range("D10").Value = "'0.01412"
$ws.Range("D12").Value = "'0.08693"
$ws.Range("D13").Value = "'0.03716"
$ws.Range("D14").Value = "'0.03203"
$ws.Range("D15").Value = "'0.09183"
$ws.Range("D16").Value = "'3.701"
$ws.Range("D17").Value = "'0.001649"
$ws.Range("D18").Value = "'0.04746"
$ws.Range("D19").Value = "'0.006197"
$ws.Range("D20").Value = "'0.006288"
$ws.Range("D23").Value = "'3.785"
$ws.Range("D40").Value = "'0.04786"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007157"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.004505"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1115"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.01168"
$ws.Range("D45").Value = "'0.00006939"
$ws.Range("D48").Value = "'0.004367"
